# Big update for 02Jan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 30-44 with revised values ---

# Row 30
$ws.Cells.Item(30, 5).Value = 0.15109999999999998
$ws.Cells.Item(30, 6).Value = 0.13439999999999999

# Row 32
$ws.Cells.Item(32, 3).Value = 0.2132
$ws.Cells.Item(32, 4).Value = 0.17100000000000001
$ws.Cells.Item(32, 6).Value = 0.12609999999999999

# Row 33
$ws.Cells.Item(33, 2).Value = 0.1787
$ws.Cells.Item(33, 4).Value = 0.1691

# Row 34
$ws.Cells.Item(34, 3).Value = 0.18129999999999999
$ws.Cells.Item(34, 4).Value = 0.1714

# Row 35
$ws.Cells.Item(35, 3).Value = 0.17929999999999999
$ws.Cells.Item(35, 4).Value = 0.17859999999999998

# Row 36
$ws.Cells.Item(36, 5).Value = 0.14949999999999999
$ws.Cells.Item(36, 6).Value = 0.13869999999999999
$ws.Cells.Item(36, 7).Value = 0.0925

# Row 38
$ws.Cells.Item(38, 6).Value = 0.14050000000000001

# Row 39
$ws.Cells.Item(39, 5).Value = 0.13400000000000001
$ws.Cells.Item(39, 6).Value = 0.1424
$ws.Cells.Item(39, 7).Value = 0.0966

# Row 40
$ws.Cells.Item(40, 2).Value = 0.1578
$ws.Cells.Item(40, 3).Value = 0.1923
$ws.Cells.Item(40, 4).Value = 0.16320000000000001
$ws.Cells.Item(40, 5).Value = 0.14069999999999999
$ws.Cells.Item(40, 7).Value = 0.099600000000000008
$ws.Cells.Item(40, 8).Value = 0.0557

# Row 41
$ws.Cells.Item(41, 2).Value = 0.16649999999999998
$ws.Cells.Item(41, 4).Value = 0.16760000000000003
$ws.Cells.Item(41, 6).Value = 0.1472
$ws.Cells.Item(41, 7).Value = 0.095399999999999985

# Row 42
$ws.Cells.Item(42, 2).Value = 0.16889999999999999
$ws.Cells.Item(42, 3).Value = 0.19159999999999999
$ws.Cells.Item(42, 6).Value = 0.14829999999999999
$ws.Cells.Item(42, 8).Value = 0.050799999999999998
$ws.Cells.Item(42, 9).Value = 0.0407

# Row 43
$ws.Cells.Item(43, 2).Value = 0.1648
$ws.Cells.Item(43, 3).Value = 0.18329999999999999
$ws.Cells.Item(43, 4).Value = 0.15789999999999998
$ws.Cells.Item(43, 5).Value = 0.1386
$ws.Cells.Item(43, 6).Value = 0.1522
$ws.Cells.Item(43, 7).Value = 0.1048
$ws.Cells.Item(43, 8).Value = 0.054299999999999994
$ws.Cells.Item(43, 9).Value = 0.043299999999999998
$ws.Cells.Item(43, 10).Value = 0.00080000000000000004

# Row 44
$ws.Cells.Item(44, 2).Value = 0.16250000000000001
$ws.Cells.Item(44, 3).Value = 0.1709
$ws.Cells.Item(44, 4).Value = 0.15049999999999999
$ws.Cells.Item(44, 5).Value = 0.13739999999999999
$ws.Cells.Item(44, 6).Value = 0.1545
$ws.Cells.Item(44, 7).Value = 0.114
$ws.Cells.Item(44, 8).Value = 0.059699999999999996
$ws.Cells.Item(44, 9).Value = 0.050199999999999995
$ws.Cells.Item(44, 10).Value = 0.00029999999999999997

# --- Append new row 45 ---
$ws.Cells.Item(45, 1).Value = 44192
$ws.Cells.Item(45, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(45, 2).Value = 0.16800000000000001
$ws.Cells.Item(45, 3).Value = 0.17530000000000001
$ws.Cells.Item(45, 4).Value = 0.16200000000000001
$ws.Cells.Item(45, 5).Value = 0.14369999999999999
$ws.Cells.Item(45, 6).Value = 0.15
$ws.Cells.Item(45, 7).Value = 0.1012
$ws.Cells.Item(45, 8).Value = 0.057300000000000004
$ws.Cells.Item(45, 9).Value = 0.041900000000000007
$ws.Cells.Item(45, 10).Value = 0.00059999999999999995

# --- Update selection to reflect saved view state ---
$ws.Range("M12").Select()
